$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Annotate the effect-column header choice next to the benefit legend.
$ws.Range("M13").Value = "effect=HR"

# Relabel the HR/SE header cells to effect/se (the underlying shared-string
# pool entries for the old "HR" and "SE" labels become unused and are
# replaced by "effect", and "se").
$ws.Range("G1").Value = "effect"
$ws.Range("J1").Value = "se"

# Add the new "se" column: se = (CIH - CIL) / 1.96^2, computed per-row.
# Enter row 2 and row 3 individually (matching how the author built this up),
# then fill the rest of the column (rows 4-17) as one contiguous block so it
# becomes a single shared formula group.
$ws.Range("J2").Formula = "=(I2-H2)/(1.96^2)"
$ws.Range("J3").Formula = "=(I3-H3)/(1.96^2)"
$ws.Range("J4:J17").Formula = "=(I4-H4)/(1.96^2)"

# Move the active selection to K17.
$ws.Range("K17").Select()
